# feat: Completed add_invoice method
#
# Adds a new invoice row to the "June" calendar sheet: row 4 gets the
# client label "EIHAB - 27 Tiff IRA" in column A and the invoice code
# "170710-SK" across columns C, D, J and K (mirroring row 3's layout).
# The existing row 3 entries that previously read "S K" are updated to
# the new invoice code "170710-SK" as well (same shared-string cell
# value, now corrected/completed).
# Finally, the active selection on the "Sheet2" tab moves to G22.

$wb = $excel.ActiveWorkbook

$juneSheet = $wb.Worksheets.Item("June")

# Update the pre-existing "S K" entries in row 3 to the completed
# invoice code "170710-SK".
$juneSheet.Range("C3").Value = "170710-SK"
$juneSheet.Range("D3").Value = "170710-SK"
$juneSheet.Range("J3").Value = "170710-SK"
$juneSheet.Range("K3").Value = "170710-SK"

# Add the new invoice row.
$juneSheet.Range("A4").Value = "EIHAB - 27 Tiff IRA"
$juneSheet.Range("C4").Value = "170710-SK"
$juneSheet.Range("D4").Value = "170710-SK"
$juneSheet.Range("J4").Value = "170710-SK"
$juneSheet.Range("K4").Value = "170710-SK"

# Move the active selection on the currently-selected tab ("Sheet2") to G22.
$sheet2 = $wb.Worksheets.Item("Sheet2")
$sheet2.Activate()
$sheet2.Range("G22").Select()
